# Apply targeted cell value toggles in column A as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
@{Row=32; Val=0}; @{Row=70; Val=0}; @{Row=202; Val=1}; @{Row=211; Val=0}; @{Row=309; Val=1}; @{Row=325; Val=1}; @{Row=462; Val=0}; @{Row=550; Val=0}; @{Row=614; Val=1}; @{Row=640; Val=1}; @{Row=715; Val=1}; @{Row=717; Val=0}; @{Row=735; Val=0}; @{Row=737; Val=0}; @{Row=741; Val=1}; @{Row=747; Val=1}; @{Row=750; Val=1}; @{Row=758; Val=0}; @{Row=766; Val=0}; @{Row=779; Val=1}; @{Row=781; Val=1}; @{Row=798; Val=1}; @{Row=811; Val=1}; @{Row=817; Val=0}; @{Row=826; Val=1}; @{Row=830; Val=1}; @{Row=842; Val=0}; @{Row=845; Val=0}; @{Row=857; Val=1}; @{Row=859; Val=1}; @{Row=860; Val=0}; @{Row=865; Val=0}; @{Row=868; Val=0}; @{Row=870; Val=0}; @{Row=875; Val=0}; @{Row=904; Val=0}; @{Row=910; Val=0}; @{Row=916; Val=0}; @{Row=950; Val=0}; @{Row=953; Val=0}; @{Row=1067; Val=1}; @{Row=1078; Val=0}; @{Row=1170; Val=0}; @{Row=1260; Val=1}; @{Row=1280; Val=1}; @{Row=1282; Val=1}; @{Row=1286; Val=1}; @{Row=1307; Val=1}; @{Row=1313; Val=1}; @{Row=1324; Val=1}; @{Row=1358; Val=0}; @{Row=1362; Val=1}; @{Row=1382; Val=0}; @{Row=1390; Val=1}; @{Row=1423; Val=1}; @{Row=1426; Val=0}; @{Row=1441; Val=0}; @{Row=1476; Val=1}; @{Row=1527; Val=1}; @{Row=1627; Val=0}; @{Row=1628; Val=0}; @{Row=1629; Val=0}; @{Row=1630; Val=0}; @{Row=1631; Val=0}; @{Row=1632; Val=0}; @{Row=1633; Val=0}; @{Row=1634; Val=0}; @{Row=1635; Val=0}; @{Row=1637; Val=1}; @{Row=1641; Val=0}; @{Row=1644; Val=0}; @{Row=1647; Val=1}; @{Row=1648; Val=0}; @{Row=1651; Val=0}; @{Row=1661; Val=0}; @{Row=1665; Val=0}; @{Row=1667; Val=1}; @{Row=1669; Val=0}; @{Row=1675; Val=1}; @{Row=1676; Val=1}; @{Row=1677; Val=0}; @{Row=1680; Val=0}; @{Row=1685; Val=0}; @{Row=1691; Val=1}; @{Row=1693; Val=1}; @{Row=1695; Val=1}; @{Row=1698; Val=1}; @{Row=1700; Val=0}; @{Row=1701; Val=1}; @{Row=1704; Val=1}; @{Row=1705; Val=0}; @{Row=1708; Val=0}; @{Row=1709; Val=1}; @{Row=1710; Val=0}; @{Row=1721; Val=1}; @{Row=1725; Val=0}; @{Row=1737; Val=1}; @{Row=1740; Val=0}; @{Row=1741; Val=1}; @{Row=1743; Val=0}; @{Row=1748; Val=1}; @{Row=1753; Val=1}; @{Row=1760; Val=1}; @{Row=1762; Val=1}; @{Row=1763; Val=1}; @{Row=1769; Val=1}; @{Row=1773; Val=1}; @{Row=1774; Val=1}; @{Row=1775; Val=1}; @{Row=1789; Val=1}; @{Row=1792; Val=0}; @{Row=1793; Val=0}; @{Row=1798; Val=0}; @{Row=1799; Val=0}; @{Row=1801; Val=0}
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg.Row, 1).Value = $chg.Val
}
